# Add a new generated "line-graph" picture to slide 1, mirroring the two
# existing generated-shape pictures already on that slide.
#
# The cleanest way to reproduce the target <p:pic> block (custom cNvPr
# extLst marking it as a generated asset, the a14:useLocalDpi blip
# extension, picLocks, etc.) through the PowerPoint object model is to
# duplicate one of the existing generated pictures -- Duplicate() carries
# all of that OOXML along with it (and mints a fresh image relationship) --
# and then fix up the few properties that differ: the shape name and its
# position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find an existing "generated asset" picture shape to clone.
$source = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Type -eq 13 -and $candidate.Name -eq "My Shape") {
        $source = $candidate
    }
}
if ($source -eq $null) {
    $source = $s.Shapes.Item($s.Shapes.Count)
}

$newPic = $source.Duplicate()
$newPic.Name = "Generated Shape"
$newPic.Left = 0
$newPic.Top = 0
